$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 22
$ws.Range("B2").Value = 112

$ws.Range("A3").Value = 21
$ws.Range("B3").Value = 103

$ws.Range("A4").Value = 11
$ws.Range("B4").Value = 69

$ws.Range("A5").Value = 12
$ws.Range("B5").Value = 14
